$wb = $excel.ActiveWorkbook

# Sheet 1 - year 2000
$ws = $wb.Worksheets.Item(1)
$ws.Range("C2").Value = 5943366
$ws.Range("D2").Value = 440.97274780273438
$ws.Range("E2").Value = 0.98911607265472412
$ws.Range("C3").Value = 5943366
$ws.Range("D3").Value = 440.97274780273438
$ws.Range("E3").Value = 0.99111831188201904
$ws.Range("C4").Value = 5943366
$ws.Range("D4").Value = 440.97274780273438
$ws.Range("E4").Value = 0.99243158102035522
$ws.Range("C5").Value = 5943366
$ws.Range("D5").Value = 440.97274780273438
$ws.Range("E5").Value = 0.99577999114990234
$ws.Range("C6").Value = 5943366
$ws.Range("D6").Value = 440.97274780273438
$ws.Range("E6").Value = 0.99758821725845337
$ws.Range("C7").Value = 5943366
$ws.Range("D7").Value = 440.97274780273438
$ws.Range("E7").Value = 0.99874430894851685
$ws.Range("C8").Value = 5943366
$ws.Range("D8").Value = 440.97274780273438
$ws.Range("E8").Value = 0.99938875436782837
$ws.Range("C9").Value = 5943366
$ws.Range("D9").Value = 440.97274780273438
$ws.Range("E9").Value = 0.999519944190979
$ws.Range("C10").Value = 5943366
$ws.Range("D10").Value = 440.97274780273438
$ws.Range("E10").Value = 0.99987363815307617
$ws.Range("C11").Value = 5943366
$ws.Range("D11").Value = 440.97274780273438
$ws.Range("E11").Value = 0.99994730949401855

# Sheet 2 - year 2001
$ws = $wb.Worksheets.Item(2)
$ws.Range("C2").Value = 5968060
$ws.Range("D2").Value = 467.48733520507812
$ws.Range("E2").Value = 0.98682636022567749
$ws.Range("C3").Value = 5968060
$ws.Range("D3").Value = 467.48733520507812
$ws.Range("E3").Value = 0.9900592565536499
$ws.Range("C4").Value = 5968060
$ws.Range("D4").Value = 467.48733520507812
$ws.Range("E4").Value = 0.99182802438735962
$ws.Range("C5").Value = 5968060
$ws.Range("D5").Value = 467.48733520507812
$ws.Range("E5").Value = 0.99554228782653809
$ws.Range("C6").Value = 5968060
$ws.Range("D6").Value = 467.48733520507812
$ws.Range("E6").Value = 0.99751609563827515
$ws.Range("C7").Value = 5968060
$ws.Range("D7").Value = 467.48733520507812
$ws.Range("E7").Value = 0.99869823455810547
$ws.Range("C8").Value = 5968060
$ws.Range("D8").Value = 467.48733520507812
$ws.Range("E8").Value = 0.99934905767440796
$ws.Range("C9").Value = 5968060
$ws.Range("D9").Value = 467.48733520507812
$ws.Range("E9").Value = 0.99948173761367798
$ws.Range("C10").Value = 5968060
$ws.Range("D10").Value = 467.48733520507812
$ws.Range("E10").Value = 0.99986380338668823
$ws.Range("C11").Value = 5968060
$ws.Range("D11").Value = 467.48733520507812
$ws.Range("E11").Value = 0.99994772672653198

# Sheet 3 - year 2002
$ws = $wb.Worksheets.Item(3)
$ws.Range("C2").Value = 5986631
$ws.Range("D2").Value = 521.77850341796875
$ws.Range("E2").Value = 0.98340803384780884
$ws.Range("C3").Value = 5986631
$ws.Range("D3").Value = 521.77850341796875
$ws.Range("E3").Value = 0.98825430870056152
$ws.Range("C4").Value = 5986631
$ws.Range("D4").Value = 521.77850341796875
$ws.Range("E4").Value = 0.99050968885421753
$ws.Range("C5").Value = 5986631
$ws.Range("D5").Value = 521.77850341796875
$ws.Range("E5").Value = 0.99504613876342773
$ws.Range("C6").Value = 5986631
$ws.Range("D6").Value = 521.77850341796875
$ws.Range("E6").Value = 0.99727374315261841
$ws.Range("C7").Value = 5986631
$ws.Range("D7").Value = 521.77850341796875
$ws.Range("E7").Value = 0.99856561422348022
$ws.Range("C8").Value = 5986631
$ws.Range("D8").Value = 521.77850341796875
$ws.Range("E8").Value = 0.99928122758865356
$ws.Range("C9").Value = 5986631
$ws.Range("D9").Value = 521.77850341796875
$ws.Range("E9").Value = 0.99943053722381592
$ws.Range("C10").Value = 5986631
$ws.Range("D10").Value = 521.77850341796875
$ws.Range("E10").Value = 0.99984651803970337
$ws.Range("C11").Value = 5986631
$ws.Range("D11").Value = 521.77850341796875
$ws.Range("E11").Value = 0.99993818998336792

# Sheet 4 - year 2003
$ws = $wb.Worksheets.Item(4)
$ws.Range("C2").Value = 5998599
$ws.Range("D2").Value = 569.46942138671875
$ws.Range("E2").Value = 0.9804912805557251
$ws.Range("C3").Value = 5998599
$ws.Range("D3").Value = 569.46942138671875
$ws.Range("E3").Value = 0.98666107654571533
$ws.Range("C4").Value = 5998599
$ws.Range("D4").Value = 569.46942138671875
$ws.Range("E4").Value = 0.98936867713928223
$ws.Range("C5").Value = 5998599
$ws.Range("D5").Value = 569.46942138671875
$ws.Range("E5").Value = 0.99461174011230469
$ws.Range("C6").Value = 5998599
$ws.Range("D6").Value = 569.46942138671875
$ws.Range("E6").Value = 0.99706530570983887
$ws.Range("C7").Value = 5998599
$ws.Range("D7").Value = 569.46942138671875
$ws.Range("E7").Value = 0.99847781658172607
$ws.Range("C8").Value = 5998599
$ws.Range("D8").Value = 569.46942138671875
$ws.Range("E8").Value = 0.99922245740890503
$ws.Range("C9").Value = 5998599
$ws.Range("D9").Value = 569.46942138671875
$ws.Range("E9").Value = 0.99938619136810303
$ws.Range("C10").Value = 5998599
$ws.Range("D10").Value = 569.46942138671875
$ws.Range("E10").Value = 0.99983680248260498
$ws.Range("C11").Value = 5998599
$ws.Range("D11").Value = 569.46942138671875
$ws.Range("E11").Value = 0.99992901086807251

# Sheet 5 - year 2004
$ws = $wb.Worksheets.Item(5)
$ws.Range("C2").Value = 6004671
$ws.Range("D2").Value = 631.49896240234375
$ws.Range("E2").Value = 0.97755664587020874
$ws.Range("C3").Value = 6004671
$ws.Range("D3").Value = 631.49896240234375
$ws.Range("E3").Value = 0.98512041568756104
$ws.Range("C4").Value = 6004671
$ws.Range("D4").Value = 631.49896240234375
$ws.Range("E4").Value = 0.9882882833480835
$ws.Range("C5").Value = 6004671
$ws.Range("D5").Value = 631.49896240234375
$ws.Range("E5").Value = 0.99425148963928223
$ws.Range("C6").Value = 6004671
$ws.Range("D6").Value = 631.49896240234375
$ws.Range("E6").Value = 0.99689441919326782
$ws.Range("C7").Value = 6004671
$ws.Range("D7").Value = 631.49896240234375
$ws.Range("E7").Value = 0.99837595224380493
$ws.Range("C8").Value = 6004671
$ws.Range("D8").Value = 631.49896240234375
$ws.Range("E8").Value = 0.99917048215866089
$ws.Range("C9").Value = 6004671
$ws.Range("D9").Value = 631.49896240234375
$ws.Range("E9").Value = 0.99933171272277832
$ws.Range("C10").Value = 6004671
$ws.Range("D10").Value = 631.49896240234375
$ws.Range("E10").Value = 0.99981915950775146
$ws.Range("C11").Value = 6004671
$ws.Range("D11").Value = 631.49896240234375
$ws.Range("E11").Value = 0.99992471933364868

# Sheet 6 - year 2005
$ws = $wb.Worksheets.Item(6)
$ws.Range("C2").Value = 6005578
$ws.Range("D2").Value = 666.664794921875
$ws.Range("E2").Value = 0.97234135866165161
$ws.Range("C3").Value = 6005578
$ws.Range("D3").Value = 666.664794921875
$ws.Range("E3").Value = 0.98221355676651001
$ws.Range("C4").Value = 6005578
$ws.Range("D4").Value = 666.664794921875
$ws.Range("E4").Value = 0.98630636930465698
$ws.Range("C5").Value = 6005578
$ws.Range("D5").Value = 666.664794921875
$ws.Range("E5").Value = 0.99362939596176147
$ws.Range("C6").Value = 6005578
$ws.Range("D6").Value = 666.664794921875
$ws.Range("E6").Value = 0.99667274951934814
$ws.Range("C7").Value = 6005578
$ws.Range("D7").Value = 666.664794921875
$ws.Range("E7").Value = 0.99830156564712524
$ws.Range("C8").Value = 6005578
$ws.Range("D8").Value = 666.664794921875
$ws.Range("E8").Value = 0.99916726350784302
$ws.Range("C9").Value = 6005578
$ws.Range("D9").Value = 666.664794921875
$ws.Range("E9").Value = 0.9993441104888916
$ws.Range("C10").Value = 6005578
$ws.Range("D10").Value = 666.664794921875
$ws.Range("E10").Value = 0.99982547760009766
$ws.Range("C11").Value = 6005578
$ws.Range("D11").Value = 666.664794921875
$ws.Range("E11").Value = 0.99992555379867554

# Sheet 7 - year 2006
$ws = $wb.Worksheets.Item(7)
$ws.Range("C2").Value = 6002319
$ws.Range("D2").Value = 756.66986083984375
$ws.Range("E2").Value = 0.96931672096252441
$ws.Range("C3").Value = 6002319
$ws.Range("D3").Value = 756.66986083984375
$ws.Range("E3").Value = 0.9805225133895874
$ws.Range("C4").Value = 6002319
$ws.Range("D4").Value = 756.66986083984375
$ws.Range("E4").Value = 0.985107421875
$ws.Range("C5").Value = 6002319
$ws.Range("D5").Value = 756.66986083984375
$ws.Range("E5").Value = 0.99311733245849609
$ws.Range("C6").Value = 6002319
$ws.Range("D6").Value = 756.66986083984375
$ws.Range("E6").Value = 0.99636906385421753
$ws.Range("C7").Value = 6002319
$ws.Range("D7").Value = 756.66986083984375
$ws.Range("E7").Value = 0.99814540147781372
$ws.Range("C8").Value = 6002319
$ws.Range("D8").Value = 756.66986083984375
$ws.Range("E8").Value = 0.99907600879669189
$ws.Range("C9").Value = 6002319
$ws.Range("D9").Value = 756.66986083984375
$ws.Range("E9").Value = 0.99926745891571045
$ws.Range("C10").Value = 6002319
$ws.Range("D10").Value = 756.66986083984375
$ws.Range("E10").Value = 0.99980926513671875
$ws.Range("C11").Value = 6002319
$ws.Range("D11").Value = 756.66986083984375
$ws.Range("E11").Value = 0.99991768598556519

# Sheet 8 - year 2007
$ws = $wb.Worksheets.Item(8)
$ws.Range("C2").Value = 6009824
$ws.Range("D2").Value = 787.24798583984375
$ws.Range("E2").Value = 0.96578752994537354
$ws.Range("C3").Value = 6009824
$ws.Range("D3").Value = 787.24798583984375
$ws.Range("E3").Value = 0.97900736331939697
$ws.Range("C4").Value = 6009824
$ws.Range("D4").Value = 787.24798583984375
$ws.Range("E4").Value = 0.98426508903503418
$ws.Range("C5").Value = 6009824
$ws.Range("D5").Value = 787.24798583984375
$ws.Range("E5").Value = 0.99283558130264282
$ws.Range("C6").Value = 6009824
$ws.Range("D6").Value = 787.24798583984375
$ws.Range("E6").Value = 0.99623751640319824
$ws.Range("C7").Value = 6009824
$ws.Range("D7").Value = 787.24798583984375
$ws.Range("E7").Value = 0.99808114767074585
$ws.Range("C8").Value = 6009824
$ws.Range("D8").Value = 787.24798583984375
$ws.Range("E8").Value = 0.99904769659042358
$ws.Range("C9").Value = 6009824
$ws.Range("D9").Value = 787.24798583984375
$ws.Range("E9").Value = 0.99924010038375854
$ws.Range("C10").Value = 6009824
$ws.Range("D10").Value = 787.24798583984375
$ws.Range("E10").Value = 0.99979919195175171
$ws.Range("C11").Value = 6009824
$ws.Range("D11").Value = 787.24798583984375
$ws.Range("E11").Value = 0.99991416931152344

# Sheet 9 - year 2009
$ws = $wb.Worksheets.Item(9)
$ws.Range("C2").Value = 6048279
$ws.Range("D2").Value = 766.649169921875
$ws.Range("E2").Value = 0.96388942003250122
$ws.Range("C3").Value = 6048279
$ws.Range("D3").Value = 766.649169921875
$ws.Range("E3").Value = 0.97746050357818604
$ws.Range("C4").Value = 6048279
$ws.Range("D4").Value = 766.649169921875
$ws.Range("E4").Value = 0.98328202962875366
$ws.Range("C5").Value = 6048279
$ws.Range("D5").Value = 766.649169921875
$ws.Range("E5").Value = 0.99287348985671997
$ws.Range("C6").Value = 6048279
$ws.Range("D6").Value = 766.649169921875
$ws.Range("E6").Value = 0.99633848667144775
$ws.Range("C7").Value = 6048279
$ws.Range("D7").Value = 766.649169921875
$ws.Range("E7").Value = 0.99817103147506714
$ws.Range("C8").Value = 6048279
$ws.Range("D8").Value = 766.649169921875
$ws.Range("E8").Value = 0.99911099672317505
$ws.Range("C9").Value = 6048279
$ws.Range("D9").Value = 766.649169921875
$ws.Range("E9").Value = 0.99929565191268921
$ws.Range("C10").Value = 6048279
$ws.Range("D10").Value = 766.649169921875
$ws.Range("E10").Value = 0.99981236457824707
$ws.Range("C11").Value = 6048279
$ws.Range("D11").Value = 766.649169921875
$ws.Range("E11").Value = 0.99991416931152344

# Sheet 10 - year 2010
$ws = $wb.Worksheets.Item(10)
$ws.Range("C2").Value = 6068249
$ws.Range("D2").Value = 819.95513916015625
$ws.Range("E2").Value = 0.96091145277023315
$ws.Range("C3").Value = 6068249
$ws.Range("D3").Value = 819.95513916015625
$ws.Range("E3").Value = 0.97567236423492432
$ws.Range("C4").Value = 6068249
$ws.Range("D4").Value = 819.95513916015625
$ws.Range("E4").Value = 0.98242807388305664
$ws.Range("C5").Value = 6068249
$ws.Range("D5").Value = 819.95513916015625
$ws.Range("E5").Value = 0.99243569374084473
$ws.Range("C6").Value = 6068249
$ws.Range("D6").Value = 819.95513916015625
$ws.Range("E6").Value = 0.996124267578125
$ws.Range("C7").Value = 6068249
$ws.Range("D7").Value = 819.95513916015625
$ws.Range("E7").Value = 0.99809056520462036
$ws.Range("C8").Value = 6068249
$ws.Range("D8").Value = 819.95513916015625
$ws.Range("E8").Value = 0.99906927347183228
$ws.Range("C9").Value = 6068249
$ws.Range("D9").Value = 819.95513916015625
$ws.Range("E9").Value = 0.99926435947418213
$ws.Range("C10").Value = 6068249
$ws.Range("D10").Value = 819.95513916015625
$ws.Range("E10").Value = 0.99980288743972778
$ws.Range("C11").Value = 6068249
$ws.Range("D11").Value = 819.95513916015625
$ws.Range("E11").Value = 0.9999118447303772

# Sheet 11 - year 2012
$ws = $wb.Worksheets.Item(11)
$ws.Range("C2").Value = 6113975
$ws.Range("D2").Value = 945.66876220703125
$ws.Range("E2").Value = 0.964607834815979
$ws.Range("C3").Value = 6113975
$ws.Range("D3").Value = 945.66876220703125
$ws.Range("E3").Value = 0.97646230459213257
$ws.Range("C4").Value = 6113975
$ws.Range("D4").Value = 945.66876220703125
$ws.Range("E4").Value = 0.98300808668136597
$ws.Range("C5").Value = 6113975
$ws.Range("D5").Value = 945.66876220703125
$ws.Range("E5").Value = 0.99190658330917358
$ws.Range("C6").Value = 6113975
$ws.Range("D6").Value = 945.66876220703125
$ws.Range("E6").Value = 0.99575608968734741
$ws.Range("C7").Value = 6113975
$ws.Range("D7").Value = 945.66876220703125
$ws.Range("E7").Value = 0.99787521362304688
$ws.Range("C8").Value = 6113975
$ws.Range("D8").Value = 945.66876220703125
$ws.Range("E8").Value = 0.99894964694976807
$ws.Range("C9").Value = 6113975
$ws.Range("D9").Value = 945.66876220703125
$ws.Range("E9").Value = 0.99917453527450562
$ws.Range("C10").Value = 6113975
$ws.Range("D10").Value = 945.66876220703125
$ws.Range("E10").Value = 0.99978232383728027
$ws.Range("C11").Value = 6113975
$ws.Range("D11").Value = 945.66876220703125
$ws.Range("E11").Value = 0.99989646673202515

# Sheet 12 - year 2013
$ws = $wb.Worksheets.Item(12)
$ws.Range("C2").Value = 6138839
$ws.Range("D2").Value = 926.16314697265625
$ws.Range("E2").Value = 0.96369034051895142
$ws.Range("C3").Value = 6138839
$ws.Range("D3").Value = 926.16314697265625
$ws.Range("E3").Value = 0.97551459074020386
$ws.Range("C4").Value = 6138839
$ws.Range("D4").Value = 926.16314697265625
$ws.Range("E4").Value = 0.98259425163269043
$ws.Range("C5").Value = 6138839
$ws.Range("D5").Value = 926.16314697265625
$ws.Range("E5").Value = 0.99171930551528931
$ws.Range("C6").Value = 6138839
$ws.Range("D6").Value = 926.16314697265625
$ws.Range("E6").Value = 0.99567180871963501
$ws.Range("C7").Value = 6138839
$ws.Range("D7").Value = 926.16314697265625
$ws.Range("E7").Value = 0.99785661697387695
$ws.Range("C8").Value = 6138839
$ws.Range("D8").Value = 926.16314697265625
$ws.Range("E8").Value = 0.998931884765625
$ws.Range("C9").Value = 6138839
$ws.Range("D9").Value = 926.16314697265625
$ws.Range("E9").Value = 0.99915796518325806
$ws.Range("C10").Value = 6138839
$ws.Range("D10").Value = 926.16314697265625
$ws.Range("E10").Value = 0.99978578090667725
$ws.Range("C11").Value = 6138839
$ws.Range("D11").Value = 926.16314697265625
$ws.Range("E11").Value = 0.99989867210388184

# Sheet 13 - year 2014
$ws = $wb.Worksheets.Item(13)
$ws.Range("C2").Value = 6162955
$ws.Range("D2").Value = 921.06500244140625
$ws.Range("E2").Value = 0.96198904514312744
$ws.Range("C3").Value = 6162955
$ws.Range("D3").Value = 921.06500244140625
$ws.Range("E3").Value = 0.97436702251434326
$ws.Range("C4").Value = 6162955
$ws.Range("D4").Value = 921.06500244140625
$ws.Range("E4").Value = 0.98203831911087036
$ws.Range("C5").Value = 6162955
$ws.Range("D5").Value = 921.06500244140625
$ws.Range("E5").Value = 0.9916192889213562
$ws.Range("C6").Value = 6162955
$ws.Range("D6").Value = 921.06500244140625
$ws.Range("E6").Value = 0.99566459655761719
$ws.Range("C7").Value = 6162955
$ws.Range("D7").Value = 921.06500244140625
$ws.Range("E7").Value = 0.99788331985473633
$ws.Range("C8").Value = 6162955
$ws.Range("D8").Value = 921.06500244140625
$ws.Range("E8").Value = 0.99895393848419189
$ws.Range("C9").Value = 6162955
$ws.Range("D9").Value = 921.06500244140625
$ws.Range("E9").Value = 0.99917310476303101
$ws.Range("C10").Value = 6162955
$ws.Range("D10").Value = 921.06500244140625
$ws.Range("E10").Value = 0.99979490041732788
$ws.Range("C11").Value = 6162955
$ws.Range("D11").Value = 921.06500244140625
$ws.Range("E11").Value = 0.99989938735961914

# Sheet 14 - year 2015
$ws = $wb.Worksheets.Item(14)
$ws.Range("C2").Value = 6183676
$ws.Range("D2").Value = 906.8468017578125
$ws.Range("E2").Value = 0.96065366268157959
$ws.Range("C3").Value = 6183676
$ws.Range("D3").Value = 906.8468017578125
$ws.Range("E3").Value = 0.97343325614929199
$ws.Range("C4").Value = 6183676
$ws.Range("D4").Value = 906.8468017578125
$ws.Range("E4").Value = 0.98140639066696167
$ws.Range("C5").Value = 6183676
$ws.Range("D5").Value = 906.8468017578125
$ws.Range("E5").Value = 0.9913485050201416
$ws.Range("C6").Value = 6183676
$ws.Range("D6").Value = 906.8468017578125
$ws.Range("E6").Value = 0.99553745985031128
$ws.Range("C7").Value = 6183676
$ws.Range("D7").Value = 906.8468017578125
$ws.Range("E7").Value = 0.99781310558319092
$ws.Range("C8").Value = 6183676
$ws.Range("D8").Value = 906.8468017578125
$ws.Range("E8").Value = 0.99893879890441895
$ws.Range("C9").Value = 6183676
$ws.Range("D9").Value = 906.8468017578125
$ws.Range("E9").Value = 0.99917668104171753
$ws.Range("C10").Value = 6183676
$ws.Range("D10").Value = 906.8468017578125
$ws.Range("E10").Value = 0.99979686737060547
$ws.Range("C11").Value = 6183676
$ws.Range("D11").Value = 906.8468017578125
$ws.Range("E11").Value = 0.99990379810333252

# Sheet 15 - year 2016
$ws = $wb.Worksheets.Item(15)
$ws.Range("C2").Value = 6200800
$ws.Range("D2").Value = 920.370849609375
$ws.Range("E2").Value = 0.96060091257095337
$ws.Range("C3").Value = 6200800
$ws.Range("D3").Value = 920.370849609375
$ws.Range("E3").Value = 0.97210055589675903
$ws.Range("C4").Value = 6200800
$ws.Range("D4").Value = 920.370849609375
$ws.Range("E4").Value = 0.98065328598022461
$ws.Range("C5").Value = 6200800
$ws.Range("D5").Value = 920.370849609375
$ws.Range("E5").Value = 0.99119096994400024
$ws.Range("C6").Value = 6200800
$ws.Range("D6").Value = 920.370849609375
$ws.Range("E6").Value = 0.99543303251266479
$ws.Range("C7").Value = 6200800
$ws.Range("D7").Value = 920.370849609375
$ws.Range("E7").Value = 0.99775612354278564
$ws.Range("C8").Value = 6200800
$ws.Range("D8").Value = 920.370849609375
$ws.Range("E8").Value = 0.998923659324646
$ws.Range("C9").Value = 6200800
$ws.Range("D9").Value = 920.370849609375
$ws.Range("E9").Value = 0.9991641640663147
$ws.Range("C10").Value = 6200800
$ws.Range("D10").Value = 920.370849609375
$ws.Range("E10").Value = 0.99979096651077271
$ws.Range("C11").Value = 6200800
$ws.Range("D11").Value = 920.370849609375
$ws.Range("E11").Value = 0.99990564584732056

# Sheet 16 - year 2017
$ws = $wb.Worksheets.Item(16)
$ws.Range("C2").Value = 6213533
$ws.Range("D2").Value = 953.46197509765625
$ws.Range("E2").Value = 0.95926582813262939
$ws.Range("C3").Value = 6213533
$ws.Range("D3").Value = 953.46197509765625
$ws.Range("E3").Value = 0.97099220752716064
$ws.Range("C4").Value = 6213533
$ws.Range("D4").Value = 953.46197509765625
$ws.Range("E4").Value = 0.97980570793151855
$ws.Range("C5").Value = 6213533
$ws.Range("D5").Value = 953.46197509765625
$ws.Range("E5").Value = 0.99094188213348389
$ws.Range("C6").Value = 6213533
$ws.Range("D6").Value = 953.46197509765625
$ws.Range("E6").Value = 0.99530571699142456
$ws.Range("C7").Value = 6213533
$ws.Range("D7").Value = 953.46197509765625
$ws.Range("E7").Value = 0.99769389629364014
$ws.Range("C8").Value = 6213533
$ws.Range("D8").Value = 953.46197509765625
$ws.Range("E8").Value = 0.99888068437576294
$ws.Range("C9").Value = 6213533
$ws.Range("D9").Value = 953.46197509765625
$ws.Range("E9").Value = 0.99912899732589722
$ws.Range("C10").Value = 6213533
$ws.Range("D10").Value = 953.46197509765625
$ws.Range("E10").Value = 0.99978399276733398
$ws.Range("C11").Value = 6213533
$ws.Range("D11").Value = 953.46197509765625
$ws.Range("E11").Value = 0.99990618228912354
